$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 20834716
$ws.Range("I33").Value = 25000642
$ws.Range("J33").Value = 5089.5
$ws.Range("K33").Value = 25000642
$ws.Range("L33").Value = 5089.5
$ws.Range("M33").Value = -25000413
$ws.Range("N33").Value = -5547.5

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7174.3335
$ws.Range("I62").Value = 5609.3
$ws.Range("J62").Value = 14999.5
$ws.Range("K62").Value = 5609.3
$ws.Range("L62").Value = 14999.5
$ws.Range("M62").Value = -4985.3
$ws.Range("N62").Value = -16247.5

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 7174.3335
$ws.Range("I65").Value = 5609.3
$ws.Range("J65").Value = 14999.5
$ws.Range("K65").Value = 28046.5
$ws.Range("L65").Value = 74997.5
$ws.Range("M65").Value = -24926.5
$ws.Range("N65").Value = -81237.5

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1482.6428
$ws.Range("I92").Value = 2276.2222
$ws.Range("K92").Value = 2276.2222
$ws.Range("M92").Value = -1028.2222

# ARM row 104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 9500.5
$ws.Range("J104").Value = 9500.5
$ws.Range("L104").Value = 9500.5
$ws.Range("N104").Value = -16488.5

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4672.625
$ws.Range("I122").Value = 4198
$ws.Range("J122").Value = 7995
$ws.Range("K122").Value = 12594
$ws.Range("L122").Value = 23985
$ws.Range("M122").Value = -10144
$ws.Range("N122").Value = -28885

# BSM row 126
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 100000
$ws.Range("J126").Value = 100000
$ws.Range("L126").Value = 100000
$ws.Range("N126").Value = -109880

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9204.333000000001
$ws.Range("I134").Value = 4440.7144
$ws.Range("K134").Value = 13322.1432
$ws.Range("M134").Value = -10787.1432

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3287.1365
$ws.Range("I122").Value = 2591.2144
$ws.Range("K122").Value = 7773.6432
$ws.Range("M122").Value = -5323.6432

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 30309964
$ws.Range("I134").Value = 1851.9333
$ws.Range("K134").Value = 5555.7999
$ws.Range("M134").Value = -3020.7999

# CUL row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 5346.3335
$ws.Range("I22").Value = 5019.5
$ws.Range("J22").Value = 6000
$ws.Range("K22").Value = 15058.5
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = -14889.5
$ws.Range("N22").Value = -18338

# CUL row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1729
$ws.Range("I26").Value = 2060.8
$ws.Range("J26").Value = 899.5
$ws.Range("K26").Value = 6182.400000000001
$ws.Range("L26").Value = 2698.5
$ws.Range("M26").Value = -5894.400000000001
$ws.Range("N26").Value = -3274.5

# CUL row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 5346.3335
$ws.Range("I27").Value = 5019.5
$ws.Range("J27").Value = 6000
$ws.Range("K27").Value = 15058.5
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = -14956.5
$ws.Range("N27").Value = -18204

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 747.3333
$ws.Range("J92").Value = 603.7143
$ws.Range("L92").Value = 1811.1429
$ws.Range("N92").Value = -4307.1429

# GSM row 4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3573.3462
$ws.Range("I97").Value = 3368.3
$ws.Range("J97").Value = 4256.8335
$ws.Range("K97").Value = 3368.3
$ws.Range("L97").Value = 4256.8335
$ws.Range("M97").Value = -2872.3
$ws.Range("N97").Value = -5248.8335

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3209.3125
$ws.Range("I102").Value = 3463.44
$ws.Range("J102").Value = 2301.7144
$ws.Range("K102").Value = 3463.44
$ws.Range("L102").Value = 2301.7144
$ws.Range("M102").Value = -1841.44
$ws.Range("N102").Value = -5545.7144

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1185
$ws.Range("I122").Value = 1074.2222
$ws.Range("J122").Value = 1517.3334
$ws.Range("K122").Value = 3222.6666
$ws.Range("L122").Value = 4552.0002
$ws.Range("M122").Value = -772.6665999999996
$ws.Range("N122").Value = -9452.0002

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8968.454
$ws.Range("I126").Value = 10001.571
$ws.Range("J126").Value = 8486.333000000001
$ws.Range("K126").Value = 30004.713
$ws.Range("L126").Value = 25458.999
$ws.Range("M126").Value = -27534.713
$ws.Range("N126").Value = -30398.999

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9006
$ws.Range("I132").Value = 5158.25
$ws.Range("K132").Value = 15474.75
$ws.Range("M132").Value = -12944.75

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2909.6667
$ws.Range("I16").Value = 3112.7144
$ws.Range("K16").Value = 3112.7144
$ws.Range("M16").Value = -2942.7144

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4955.393
$ws.Range("I40").Value = 2127.9
$ws.Range("J40").Value = 12024.125
$ws.Range("K40").Value = 2127.9
$ws.Range("L40").Value = 12024.125
$ws.Range("M40").Value = -1991.9
$ws.Range("N40").Value = -12296.125

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5110996.5
$ws.Range("I68").Value = 2344.5
$ws.Range("J68").Value = 7154457.5
$ws.Range("K68").Value = 2344.5
$ws.Range("L68").Value = 7154457.5
$ws.Range("M68").Value = -1595.5
$ws.Range("N68").Value = -7155955.5

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5110996.5
$ws.Range("I71").Value = 2344.5
$ws.Range("J71").Value = 7154457.5
$ws.Range("K71").Value = 11722.5
$ws.Range("L71").Value = 35772287.5
$ws.Range("M71").Value = -7978.5
$ws.Range("N71").Value = -35779775.5

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 16927.857
$ws.Range("I93").Value = 17888
$ws.Range("J93").Value = 15199.6
$ws.Range("K93").Value = 17888
$ws.Range("L93").Value = 15199.6
$ws.Range("M93").Value = -16640
$ws.Range("N93").Value = -17695.6

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3225
$ws.Range("I100").Value = 2445.5557
$ws.Range("J100").Value = 4227.143
$ws.Range("K100").Value = 2445.5557
$ws.Range("L100").Value = 4227.143
$ws.Range("M100").Value = -1904.5557
$ws.Range("N100").Value = -5309.143

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1184296.4
$ws.Range("I132").Value = 4423.6665
$ws.Range("K132").Value = 13270.9995
$ws.Range("M132").Value = -10740.9995

# WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3432.7778
$ws.Range("I2").Value = 3432.7778
$ws.Range("K2").Value = 3432.7778
$ws.Range("M2").Value = -3320.7778

# WVR row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 36279.355
$ws.Range("I4").Value = 208.66667
$ws.Range("J4").Value = 63332.375
$ws.Range("K4").Value = 208.66667
$ws.Range("L4").Value = 63332.375
$ws.Range("M4").Value = -95.66667000000001
$ws.Range("N4").Value = -63558.375

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3523.9429
$ws.Range("I122").Value = 2334.9473
$ws.Range("K122").Value = 7004.841899999999
$ws.Range("M122").Value = -4554.841899999999

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12609.546
$ws.Range("I132").Value = 2083.3333
$ws.Range("K132").Value = 6249.999899999999
$ws.Range("M132").Value = -3719.999899999999

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 12466.772
$ws.Range("I136").Value = 2019.9231
$ws.Range("K136").Value = 6059.7693
$ws.Range("M136").Value = -3509.7693

# WVR row 138
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 133331.67
$ws.Range("J138").Value = 133331.67
$ws.Range("L138").Value = 133331.67
$ws.Range("N138").Value = -143611.67
